$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos symbol-list refresh: update Price (D), Volume(1h) (E) and Hora (G)
# columns per row. Values are written as literal text (not numbers/percentages)
# to match the source data feed, which stores these columns as plain strings -
# hence the NumberFormat "@" (Text) trick before the write, and resetting the
# style back to "Normal" afterwards so no stray formatting is introduced.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "301.08"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.11%"
$ws.Range("E2").Style = "Normal"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "9"
$ws.Range("G2").Style = "Normal"
# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "32.78"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "4.08%"
$ws.Range("E3").Style = "Normal"
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "9"
$ws.Range("G3").Style = "Normal"
# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.935"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-2.98%"
$ws.Range("E4").Style = "Normal"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "9"
$ws.Range("G4").Style = "Normal"
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07702"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-1.95%"
$ws.Range("E5").Style = "Normal"
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "9"
$ws.Range("G5").Style = "Normal"
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.981"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-15.54%"
$ws.Range("E6").Style = "Normal"
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = "9"
$ws.Range("G6").Style = "Normal"
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "7.831"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.23%"
$ws.Range("E7").Style = "Normal"
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = "9"
$ws.Range("G7").Style = "Normal"
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.802"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.57%"
$ws.Range("E8").Style = "Normal"
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "9"
$ws.Range("G8").Style = "Normal"
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9202"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "0.25%"
$ws.Range("E9").Style = "Normal"
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "9"
$ws.Range("G9").Style = "Normal"
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1752"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "0.41%"
$ws.Range("E10").Style = "Normal"
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "9"
$ws.Range("G10").Style = "Normal"
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07760"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "2.44%"
$ws.Range("E11").Style = "Normal"
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "9"
$ws.Range("G11").Style = "Normal"
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08623"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-6.54%"
$ws.Range("E12").Style = "Normal"
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = "9"
$ws.Range("G12").Style = "Normal"
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03170"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "6.02%"
$ws.Range("E13").Style = "Normal"
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "9"
$ws.Range("G13").Style = "Normal"
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.1003"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.01%"
$ws.Range("E14").Style = "Normal"
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "9"
$ws.Range("G14").Style = "Normal"
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001516"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.34%"
$ws.Range("E15").Style = "Normal"
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = "9"
$ws.Range("G15").Style = "Normal"
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005910"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-0.53%"
$ws.Range("E16").Style = "Normal"
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = "9"
$ws.Range("G16").Style = "Normal"
# Row 17
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.36%"
$ws.Range("E17").Style = "Normal"
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = "9"
$ws.Range("G17").Style = "Normal"
# Row 18
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-4.19%"
$ws.Range("E18").Style = "Normal"
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = "9"
$ws.Range("G18").Style = "Normal"
# Row 19
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "2.29%"
$ws.Range("E19").Style = "Normal"
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = "9"
$ws.Range("G19").Style = "Normal"
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1328"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "3.10%"
$ws.Range("E20").Style = "Normal"
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = "9"
$ws.Range("G20").Style = "Normal"
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.309"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "6.45%"
$ws.Range("E21").Style = "Normal"
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = "9"
$ws.Range("G21").Style = "Normal"
# Row 22
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "17.24%"
$ws.Range("E22").Style = "Normal"
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = "9"
$ws.Range("G22").Style = "Normal"
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04530"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-1.83%"
$ws.Range("E23").Style = "Normal"
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = "9"
$ws.Range("G23").Style = "Normal"
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001224"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-1.83%"
$ws.Range("E24").Style = "Normal"
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = "9"
$ws.Range("G24").Style = "Normal"
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004407"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-1.43%"
$ws.Range("E25").Style = "Normal"
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = "9"
$ws.Range("G25").Style = "Normal"
# Row 26
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.26%"
$ws.Range("E26").Style = "Normal"
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = "9"
$ws.Range("G26").Style = "Normal"
# Row 27
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = "9"
$ws.Range("G27").Style = "Normal"
# Row 28
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = "9"
$ws.Range("G28").Style = "Normal"
# Row 29
$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = "9"
$ws.Range("G29").Style = "Normal"
# Row 30
$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = "9"
$ws.Range("G30").Style = "Normal"
# Row 31
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = "9"
$ws.Range("G31").Style = "Normal"
# Row 32
$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = "9"
$ws.Range("G32").Style = "Normal"
# Row 33
$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = "9"
$ws.Range("G33").Style = "Normal"
# Row 34
$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = "9"
$ws.Range("G34").Style = "Normal"
# Row 35
$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = "9"
$ws.Range("G35").Style = "Normal"
# Row 36
$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = "9"
$ws.Range("G36").Style = "Normal"
# Row 37
$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = "9"
$ws.Range("G37").Style = "Normal"
# Row 38
$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = "9"
$ws.Range("G38").Style = "Normal"
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01698"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-3.61%"
$ws.Range("E39").Style = "Normal"
$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = "9"
$ws.Range("G39").Style = "Normal"
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04678"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-0.39%"
$ws.Range("E40").Style = "Normal"
$ws.Range("G40").NumberFormat = "@"
$ws.Range("G40").Value = "9"
$ws.Range("G40").Style = "Normal"
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007593"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "5.34%"
$ws.Range("E41").Style = "Normal"
$ws.Range("G41").NumberFormat = "@"
$ws.Range("G41").Value = "9"
$ws.Range("G41").Style = "Normal"
# Row 42
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-0.88%"
$ws.Range("E42").Style = "Normal"
$ws.Range("G42").NumberFormat = "@"
$ws.Range("G42").Value = "9"
$ws.Range("G42").Style = "Normal"
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002331"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "6.57%"
$ws.Range("E43").Style = "Normal"
$ws.Range("G43").NumberFormat = "@"
$ws.Range("G43").Value = "9"
$ws.Range("G43").Style = "Normal"
# Row 44
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "17.12%"
$ws.Range("E44").Style = "Normal"
$ws.Range("G44").NumberFormat = "@"
$ws.Range("G44").Value = "9"
$ws.Range("G44").Style = "Normal"
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00006269"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "0.28%"
$ws.Range("E45").Style = "Normal"
$ws.Range("G45").NumberFormat = "@"
$ws.Range("G45").Value = "9"
$ws.Range("G45").Style = "Normal"
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000752"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.36%"
$ws.Range("E46").Style = "Normal"
$ws.Range("G46").NumberFormat = "@"
$ws.Range("G46").Value = "9"
$ws.Range("G46").Style = "Normal"
# Row 47
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "10.42%"
$ws.Range("E47").Style = "Normal"
$ws.Range("G47").NumberFormat = "@"
$ws.Range("G47").Value = "9"
$ws.Range("G47").Style = "Normal"
# Row 48
$ws.Range("G48").NumberFormat = "@"
$ws.Range("G48").Value = "9"
$ws.Range("G48").Style = "Normal"
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002105"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.36%"
$ws.Range("E49").Style = "Normal"
$ws.Range("G49").NumberFormat = "@"
$ws.Range("G49").Value = "9"
$ws.Range("G49").Style = "Normal"
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0002005"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.36%"
$ws.Range("E50").Style = "Normal"
$ws.Range("G50").NumberFormat = "@"
$ws.Range("G50").Value = "9"
$ws.Range("G50").Style = "Normal"
# Row 51
$ws.Range("G51").NumberFormat = "@"
$ws.Range("G51").Value = "9"
$ws.Range("G51").Style = "Normal"
